{"js": "// Apply the LOB1244 course-sheet update described by the diff:\n//  - Cr\u00e9ditos-aula: 2 -> 4\n//  - Carga hor\u00e1ria: 30 h -> 60 h\n//  - Ativa\u00e7\u00e3o: 01/01/2012 -> 01/01/2025\n//  - Objetivos paragraph replaced\n//  - Docente(s) list: shift + new second name\n//  - Programa resumido paragraph replaced\n//  - Programa paragraph replaced\n//  - Avalia\u00e7\u00e3o: M\u00e9todo / Crit\u00e9rio / Norma de recupera\u00e7\u00e3o values replaced\n\nasync function replaceOnce(body, searchText, newText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// Credits / workload / activation date block.\nawait replaceOnce(body, \"Cr\u00e9ditos-aula: 2\", \"Cr\u00e9ditos-aula: 4\");\nawait replaceOnce(body, \"Carga hor\u00e1ria: 30 h\", \"Carga hor\u00e1ria: 60 h\");\nawait replaceOnce(body, \"Ativa\u00e7\u00e3o: 01/01/2012\", \"Ativa\u00e7\u00e3o: 01/01/2025\");\n\n// Objetivos.\nawait replaceOnce(\n  body,\n  \"Apresentar e analisar os conceitos b\u00e1sicos de monitoramento, suas aplica\u00e7\u00f5es pr\u00e1ticas e as interfaces com os demais instrumentos de Pol\u00edtica Ambiental.\",\n  \"Apresentar e analisar as t\u00e9cnicas de monitoramento dos principais par\u00e2metros ambientais exigidos pelas legisla\u00e7\u00f5es em vigor.\"\n);\n\n// Docente(s) Respons\u00e1vel(eis) \u2014 do the second line first so the first\n// line's original text is still unique when we search for it next.\nawait replaceOnce(\n  body,\n  \"8855158 - Morun Bernardino Neto\",\n  \"7455355 - Robson da Silva Rocha\"\n);\nawait replaceOnce(\n  body,\n  \"5840938 - Marcelo Rodrigues de Holanda\",\n  \"8855158 - Morun Bernardino Neto\"\n);\n\n// Programa resumido.\nawait replaceOnce(\n  body,\n  \"Monitoramento da qualidade ambiental.\",\n  \"Conceito de monitoramento. Amostragem. T\u00e9cnicas alternativas para cada par\u00e2metro a ser monitorado. Rela\u00e7\u00e3o custo e aplicabilidade.\"\n);\n\n// Programa.\nawait replaceOnce(\n  body,\n  \"Conceitos de qualidade ambiental, polui\u00e7\u00e3o, padr\u00f5es de qualidade e de emiss\u00e3o. Conceito de monitoramento. Amostragem. Sistemas de monitoramento. \u00cdndices de qualidade. Monitoramento como parte integrante de sistema de gest\u00e3o ambiental.\",\n  \"- Conceito de monitoramento.- T\u00e9cnicas de amostragens e suas especificidades para cada par\u00e2metro.- T\u00e9cnicas e equipamentos para monitorar: carga org\u00e2nica, s\u00f3lidos, \u00edons, metais, atividade biol\u00f3gica e outros par\u00e2metros de import\u00e2ncia ambiental.- Qu\u00edmica Verde no monitoramento ambiental\"\n);\n\n// Avalia\u00e7\u00e3o block.\nawait replaceOnce(\n  body,\n  \"Aula expositiva e exerc\u00edcios dirigidos.\",\n  \"Avalia\u00e7\u00e3o baseada em provas, exerc\u00edcios, trabalhos pr\u00e1ticos e relat\u00f3rios.\"\n);\nawait replaceOnce(\n  body,\n  \"M\u00e9dia ponderada de exerc\u00edcios e provas.\",\n  \"M\u00e9dia ponderada das notas atribu\u00eddas \u00e0s provas, exerc\u00edcios e trabalhos pr\u00e1ticos e relat\u00f3rios.\"\n);\nawait replaceOnce(\n  body,\n  \"Prova \u00fanica com nota igual ou superior a 5,0.\",\n  \"1 (uma) prova de recupera\u00e7\u00e3o (R), sendo considerado aprovado se R >= 5,0.\"\n);\n", "ps1": "# Apply the LOB1244 course-sheet update described by the diff:\n#  - Cr\u00e9ditos-aula: 2 -> 4\n#  - Carga hor\u00e1ria: 30 h -> 60 h\n#  - Ativa\u00e7\u00e3o: 01/01/2012 -> 01/01/2025\n#  - Objetivos paragraph replaced\n#  - Docente(s) list: shift + new second name\n#  - Programa resumido paragraph replaced\n#  - Programa paragraph replaced\n#  - Avalia\u00e7\u00e3o: M\u00e9todo / Crit\u00e9rio / Norma de recupera\u00e7\u00e3o values replaced\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($doc, [string]$oldText, [string]$newText) {\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # Forward, Wrap=wdFindStop(0), Format=False, MatchCase=True, MatchWholeWord=False,\n    # MatchWildcards=False, MatchSoundsLike=False, MatchAllWordForms=False,\n    # Replace=wdReplaceAll(2)\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 0, $false, $newText, 2) | Out-Null\n}\n\n# Credits / workload / activation date block.\nReplace-Text $d \"Cr\u00e9ditos-aula: 2\" \"Cr\u00e9ditos-aula: 4\"\nReplace-Text $d \"Carga hor\u00e1ria: 30 h\" \"Carga hor\u00e1ria: 60 h\"\nReplace-Text $d \"Ativa\u00e7\u00e3o: 01/01/2012\" \"Ativa\u00e7\u00e3o: 01/01/2025\"\n\n# Objetivos.\nReplace-Text $d \"Apresentar e analisar os conceitos b\u00e1sicos de monitoramento, suas aplica\u00e7\u00f5es pr\u00e1ticas e as interfaces com os demais instrumentos de Pol\u00edtica Ambiental.\" \"Apresentar e analisar as t\u00e9cnicas de monitoramento dos principais par\u00e2metros ambientais exigidos pelas legisla\u00e7\u00f5es em vigor.\"\n\n# Docente(s) Respons\u00e1vel(eis) \u2014 do the second line first so the first\n# line's original text is still unique when we search for it next.\nReplace-Text $d \"8855158 - Morun Bernardino Neto\" \"7455355 - Robson da Silva Rocha\"\nReplace-Text $d \"5840938 - Marcelo Rodrigues de Holanda\" \"8855158 - Morun Bernardino Neto\"\n\n# Programa resumido.\nReplace-Text $d \"Monitoramento da qualidade ambiental.\" \"Conceito de monitoramento. Amostragem. T\u00e9cnicas alternativas para cada par\u00e2metro a ser monitorado. Rela\u00e7\u00e3o custo e aplicabilidade.\"\n\n# Programa.\nReplace-Text $d \"Conceitos de qualidade ambiental, polui\u00e7\u00e3o, padr\u00f5es de qualidade e de emiss\u00e3o. Conceito de monitoramento. Amostragem. Sistemas de monitoramento. \u00cdndices de qualidade. Monitoramento como parte integrante de sistema de gest\u00e3o ambiental.\" \"- Conceito de monitoramento.- T\u00e9cnicas de amostragens e suas especificidades para cada par\u00e2metro.- T\u00e9cnicas e equipamentos para monitorar: carga org\u00e2nica, s\u00f3lidos, \u00edons, metais, atividade biol\u00f3gica e outros par\u00e2metros de import\u00e2ncia ambiental.- Qu\u00edmica Verde no monitoramento ambiental\"\n\n# Avalia\u00e7\u00e3o block.\nReplace-Text $d \"Aula expositiva e exerc\u00edcios dirigidos.\" \"Avalia\u00e7\u00e3o baseada em provas, exerc\u00edcios, trabalhos pr\u00e1ticos e relat\u00f3rios.\"\nReplace-Text $d \"M\u00e9dia ponderada de exerc\u00edcios e provas.\" \"M\u00e9dia ponderada das notas atribu\u00eddas \u00e0s provas, exerc\u00edcios e trabalhos pr\u00e1ticos e relat\u00f3rios.\"\nReplace-Text $d \"Prova \u00fanica com nota igual ou superior a 5,0.\" \"1 (uma) prova de recupera\u00e7\u00e3o (R), sendo considerado aprovado se R >= 5,0.\"\n"}
